$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "R1"; Cell = "G2"; Value = "3929:42:41" },
    @{ Sheet = "R1"; Cell = "G3"; Value = "69:15:19" },
    @{ Sheet = "R2"; Cell = "G2"; Value = "12111:06:10" },
    @{ Sheet = "R2"; Cell = "G3"; Value = "3240:49:39" },
    @{ Sheet = "R2"; Cell = "G4"; Value = "479:01:13" },
    @{ Sheet = "R4"; Cell = "G2"; Value = "2956:55:59" },
    @{ Sheet = "R4"; Cell = "G3"; Value = "184:08:14" },
    @{ Sheet = "R4"; Cell = "G4"; Value = "72:20:39" },
    @{ Sheet = "R4"; Cell = "G5"; Value = "69:58:12" },
    @{ Sheet = "R5"; Cell = "G2"; Value = "430:54:58" },
    @{ Sheet = "R6"; Cell = "G2"; Value = "71:27:16" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
